$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I19").Value2 = 2651.1738
$ws.Range("J19").Value2 = 7260.4165
$ws.Range("K19").Value2 = 2651.1738
$ws.Range("L19").Value2 = 7260.4165
$ws.Range("M19").Value2 = -2476.1738
$ws.Range("N19").Value2 = -7610.4165
$ws.Range("H31").Value2 = 700
$ws.Range("I31").Value2 = 700
$ws.Range("J31").Value2 = 0
$ws.Range("K31").Value2 = 2100
$ws.Range("L31").Value2 = 0
$ws.Range("M31").Value2 = -1870
$ws.Range("H69").Value2 = 19419
$ws.Range("I69").Value2 = 20164.166
$ws.Range("J69").Value2 = 18922.223
$ws.Range("K69").Value2 = 60492.49800000001
$ws.Range("L69").Value2 = 56766.66900000001
$ws.Range("M69").Value2 = -59618.49800000001
$ws.Range("N69").Value2 = -58514.66900000001
$ws.Range("H70").Value2 = 2327511.2
$ws.Range("I70").Value2 = 1800
$ws.Range("J70").Value2 = 2908939
$ws.Range("K70").Value2 = 5400
$ws.Range("L70").Value2 = 8726817
$ws.Range("M70").Value2 = -5130
$ws.Range("N70").Value2 = -8727357
$ws.Range("H72").Value2 = 19419
$ws.Range("I72").Value2 = 20164.166
$ws.Range("J72").Value2 = 18922.223
$ws.Range("K72").Value2 = 181477.494
$ws.Range("L72").Value2 = 170300.007
$ws.Range("M72").Value2 = -177109.494
$ws.Range("N72").Value2 = -179036.007
$ws.Range("H73").Value2 = 2327511.2
$ws.Range("I73").Value2 = 1800
$ws.Range("J73").Value2 = 2908939
$ws.Range("K73").Value2 = 5400
$ws.Range("L73").Value2 = 8726817
$ws.Range("M73").Value2 = -4464
$ws.Range("N73").Value2 = -8728689
$ws.Range("H94").Value2 = 2966.6667
$ws.Range("I94").Value2 = 2450
$ws.Range("J94").Value2 = 4000
$ws.Range("K94").Value2 = 2450
$ws.Range("L94").Value2 = 4000
$ws.Range("M94").Value2 = -1999
$ws.Range("H103").Value2 = 1424.6
$ws.Range("I103").Value2 = 1569
$ws.Range("J103").Value2 = 1388.5
$ws.Range("K103").Value2 = 4707
$ws.Range("L103").Value2 = 4165.5
$ws.Range("M103").Value2 = -4121
$ws.Range("N103").Value2 = -5337.5
$ws.Range("H116").Value2 = 5840
$ws.Range("I116").Value2 = 5000
$ws.Range("J116").Value2 = 6400
$ws.Range("K116").Value2 = 5000
$ws.Range("L116").Value2 = 6400
$ws.Range("M116").Value2 = -1558
$ws.Range("N116").Value2 = -13284
$ws.Range("H132").Value2 = 1763.4878
$ws.Range("I132").Value2 = 1228.3125
$ws.Range("J132").Value2 = 3666.3333
$ws.Range("K132").Value2 = 3684.9375
$ws.Range("L132").Value2 = 10998.9999
$ws.Range("M132").Value2 = -1154.9375
$ws.Range("H135").Value2 = 4351.381
$ws.Range("I135").Value2 = 3748.25
$ws.Range("J135").Value2 = 5155.5557
$ws.Range("K135").Value2 = 33734.25
$ws.Range("L135").Value2 = 46400.0013
$ws.Range("M135").Value2 = -31199.25
$ws.Range("N135").Value2 = -51470.0013
$ws.Range("H138").Value2 = 5265.62
$ws.Range("I138").Value2 = 3358.3572
$ws.Range("J138").Value2 = 6007.3335
$ws.Range("K138").Value2 = 10075.0716
$ws.Range("L138").Value2 = 18022.0005
$ws.Range("M138").Value2 = -4935.071599999999
$ws.Range("N138").Value2 = -28302.0005
$ws.Range("H141").Value2 = 1985.45
$ws.Range("I141").Value2 = 1872.5555
$ws.Range("J141").Value2 = 3001.5
$ws.Range("K141").Value2 = 5617.666499999999
$ws.Range("L141").Value2 = 9004.5
$ws.Range("M141").Value2 = -437.6664999999994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 5752869
$ws.Range("I2").Value2 = 7079954
$ws.Range("J2").Value2 = 2166.3333
$ws.Range("K2").Value2 = 7079954
$ws.Range("L2").Value2 = 2166.3333
$ws.Range("M2").Value2 = -7079841
$ws.Range("N2").Value2 = -2392.3333
$ws.Range("H32").Value2 = 27806.402
$ws.Range("I32").Value2 = 25859.807
$ws.Range("J32").Value2 = 49997.6
$ws.Range("K32").Value2 = 25859.807
$ws.Range("L32").Value2 = 49997.6
$ws.Range("M32").Value2 = -25572.807
$ws.Range("N32").Value2 = -50571.6
$ws.Range("H61").Value2 = 11283.286
$ws.Range("I61").Value2 = 3439.8
$ws.Range("J61").Value2 = 15640.777
$ws.Range("K61").Value2 = 3439.8
$ws.Range("L61").Value2 = 15640.777
$ws.Range("M61").Value2 = -3227.8
$ws.Range("N61").Value2 = -16064.777
$ws.Range("H63").Value2 = 9928.799999999999
$ws.Range("I63").Value2 = 9000
$ws.Range("J63").Value2 = 9995.143
$ws.Range("K63").Value2 = 9000
$ws.Range("L63").Value2 = 9995.143
$ws.Range("M63").Value2 = -8314
$ws.Range("N63").Value2 = -11367.143
$ws.Range("H66").Value2 = 9928.799999999999
$ws.Range("I66").Value2 = 9000
$ws.Range("J66").Value2 = 9995.143
$ws.Range("K66").Value2 = 45000
$ws.Range("L66").Value2 = 49975.715
$ws.Range("M66").Value2 = -41568
$ws.Range("N66").Value2 = -56839.715
$ws.Range("H74").Value2 = 478739.34
$ws.Range("I74").Value2 = 668825.6
$ws.Range("J74").Value2 = 3523.6667
$ws.Range("K74").Value2 = 668825.6
$ws.Range("L74").Value2 = 3523.6667
$ws.Range("M74").Value2 = -667951.6
$ws.Range("N74").Value2 = -5271.6667
$ws.Range("H77").Value2 = 478739.34
$ws.Range("I77").Value2 = 668825.6
$ws.Range("J77").Value2 = 3523.6667
$ws.Range("K77").Value2 = 3344128
$ws.Range("L77").Value2 = 17618.3335
$ws.Range("M77").Value2 = -3339760
$ws.Range("N77").Value2 = -26354.3335
$ws.Range("H97").Value2 = 4744294.5
$ws.Range("I97").Value2 = 7409769
$ws.Range("J97").Value2 = 301837
$ws.Range("K97").Value2 = 7409769
$ws.Range("L97").Value2 = 301837
$ws.Range("M97").Value2 = -7409273
$ws.Range("N97").Value2 = -302829
$ws.Range("H116").Value2 = 5752869
$ws.Range("I116").Value2 = 7079954
$ws.Range("J116").Value2 = 2166.3333
$ws.Range("K116").Value2 = 7079954
$ws.Range("L116").Value2 = 2166.3333
$ws.Range("M116").Value2 = -7077660
$ws.Range("N116").Value2 = -6754.3333
$ws.Range("H122").Value2 = 3061.6
$ws.Range("I122").Value2 = 2171.1
$ws.Range("J122").Value2 = 4842.6
$ws.Range("K122").Value2 = 6513.299999999999
$ws.Range("L122").Value2 = 14527.8
$ws.Range("M122").Value2 = -4063.299999999999
$ws.Range("N122").Value2 = -19427.8
$ws.Range("H132").Value2 = 7396.6235
$ws.Range("I132").Value2 = 3140.7068
$ws.Range("J132").Value2 = 20388.37
$ws.Range("K132").Value2 = 9422.1204
$ws.Range("L132").Value2 = 61165.11
$ws.Range("M132").Value2 = -6892.1204
$ws.Range("N132").Value2 = -66225.11
$ws.Range("H136").Value2 = 11283.286
$ws.Range("I136").Value2 = 3439.8
$ws.Range("J136").Value2 = 15640.777
$ws.Range("K136").Value2 = 10319.4
$ws.Range("L136").Value2 = 46922.331
$ws.Range("M136").Value2 = -7769.400000000001
$ws.Range("N136").Value2 = -52022.331

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 5752869
$ws.Range("I3").Value2 = 7079954
$ws.Range("J3").Value2 = 2166.3333
$ws.Range("K3").Value2 = 7079954
$ws.Range("L3").Value2 = 2166.3333
$ws.Range("M3").Value2 = -7079840
$ws.Range("N3").Value2 = -2394.3333
$ws.Range("I15").Value2 = 10000
$ws.Range("J15").Value2 = 10000
$ws.Range("K15").Value2 = 10000
$ws.Range("L15").Value2 = 10000
$ws.Range("M15").Value2 = -9773
$ws.Range("N15").Value2 = -10454
$ws.Range("H19").Value2 = 15000
$ws.Range("I19").Value2 = 0
$ws.Range("J19").Value2 = 15000
$ws.Range("K19").Value2 = 0
$ws.Range("L19").Value2 = 15000
$ws.Range("M19").Value2 = -1999
$ws.Range("N19").Value2 = -15346
$ws.Range("H30").Value2 = 0
$ws.Range("I30").Value2 = 0
$ws.Range("J30").Value2 = 0
$ws.Range("K30").Value2 = 0
$ws.Range("L30").Value2 = 0
$ws.Range("M30").ClearContents()
$ws.Range("H81").Value2 = 89149.664
$ws.Range("I81").Value2 = 0
$ws.Range("J81").Value2 = 89149.664
$ws.Range("K81").Value2 = 0
$ws.Range("L81").Value2 = 89149.664
$ws.Range("N81").Value2 = -91271.664
$ws.Range("H84").Value2 = 89149.664
$ws.Range("I84").Value2 = 0
$ws.Range("J84").Value2 = 89149.664
$ws.Range("K84").Value2 = 0
$ws.Range("L84").Value2 = 267448.992
$ws.Range("N84").Value2 = -278056.992
$ws.Range("H86").Value2 = 130759.875
$ws.Range("I86").Value2 = 2205.9
$ws.Range("J86").Value2 = 345016.5
$ws.Range("K86").Value2 = 2205.9
$ws.Range("L86").Value2 = 345016.5
$ws.Range("M86").Value2 = -1082.9
$ws.Range("N86").Value2 = -347262.5
$ws.Range("H89").Value2 = 130759.875
$ws.Range("I89").Value2 = 2205.9
$ws.Range("J89").Value2 = 345016.5
$ws.Range("K89").Value2 = 11029.5
$ws.Range("L89").Value2 = 1725082.5
$ws.Range("M89").Value2 = -5413.5
$ws.Range("N89").Value2 = -1736314.5
$ws.Range("H102").Value2 = 11002.286
$ws.Range("I102").Value2 = 11002.286
$ws.Range("J102").Value2 = 0
$ws.Range("K102").Value2 = 11002.286
$ws.Range("L102").Value2 = 0
$ws.Range("M102").Value2 = -7757.286
$ws.Range("H107").Value2 = 2279.6667
$ws.Range("I107").Value2 = 1654.7142
$ws.Range("J107").Value2 = 5779.4
$ws.Range("K107").Value2 = 1654.7142
$ws.Range("L107").Value2 = 5779.4
$ws.Range("M107").Value2 = 265.2858000000001
$ws.Range("H134").Value2 = 4881.7456
$ws.Range("I134").Value2 = 4039.2554
$ws.Range("J134").Value2 = 8181.5
$ws.Range("K134").Value2 = 12117.7662
$ws.Range("L134").Value2 = 24544.5
$ws.Range("M134").Value2 = -9582.7662
$ws.Range("N134").Value2 = -29614.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 443.95
$ws.Range("I22").Value2 = 443.95
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 443.95
$ws.Range("L22").Value2 = 0
$ws.Range("M22").Value2 = -93.94999999999999
$ws.Range("H31").Value2 = 7039.206
$ws.Range("I31").Value2 = 3556.1538
$ws.Range("J31").Value2 = 9195.380999999999
$ws.Range("K31").Value2 = 3556.1538
$ws.Range("L31").Value2 = 9195.380999999999
$ws.Range("M31").Value2 = -3261.1538
$ws.Range("H34").Value2 = 7039.206
$ws.Range("I34").Value2 = 3556.1538
$ws.Range("J34").Value2 = 9195.380999999999
$ws.Range("K34").Value2 = 3556.1538
$ws.Range("L34").Value2 = 9195.380999999999
$ws.Range("M34").Value2 = -3354.1538
$ws.Range("H50").Value2 = 48247.25
$ws.Range("I50").Value2 = 50000
$ws.Range("J50").Value2 = 47663
$ws.Range("K50").Value2 = 50000
$ws.Range("L50").Value2 = 47663
$ws.Range("M50").Value2 = -49375
$ws.Range("N50").Value2 = -48913
$ws.Range("H62").Value2 = 7333.1665
$ws.Range("I62").Value2 = 7333.1665
$ws.Range("J62").Value2 = 0
$ws.Range("K62").Value2 = 7333.1665
$ws.Range("L62").Value2 = 0
$ws.Range("M62").Value2 = -6709.1665
$ws.Range("H65").Value2 = 7333.1665
$ws.Range("I65").Value2 = 7333.1665
$ws.Range("J65").Value2 = 0
$ws.Range("K65").Value2 = 36665.8325
$ws.Range("L65").Value2 = 0
$ws.Range("M65").Value2 = -33545.8325
$ws.Range("H107").Value2 = 452.1111
$ws.Range("I107").Value2 = 452.1111
$ws.Range("J107").Value2 = 0
$ws.Range("K107").Value2 = 452.1111
$ws.Range("L107").Value2 = 0
$ws.Range("M107").Value2 = 1467.8889
$ws.Range("H132").Value2 = 25587.23
$ws.Range("I132").Value2 = 3928.5
$ws.Range("J132").Value2 = 285492
$ws.Range("K132").Value2 = 11785.5
$ws.Range("L132").Value2 = 856476
$ws.Range("M132").Value2 = -9255.5
$ws.Range("H141").Value2 = 545777.2
$ws.Range("I141").Value2 = 124750
$ws.Range("J141").Value2 = 567936.5
$ws.Range("K141").Value2 = 124750
$ws.Range("L141").Value2 = 567936.5
$ws.Range("M141").Value2 = -119570
$ws.Range("N141").Value2 = -578296.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 3919.2222
$ws.Range("I2").Value2 = 50.3
$ws.Range("J2").Value2 = 8755.375
$ws.Range("K2").Value2 = 301.8
$ws.Range("L2").Value2 = 52532.25
$ws.Range("M2").Value2 = -188.8
$ws.Range("N2").Value2 = -52758.25
$ws.Range("H4").Value2 = 4345842
$ws.Range("I4").Value2 = 2826810.5
$ws.Range("J4").Value2 = 8713058
$ws.Range("K4").Value2 = 8480431.5
$ws.Range("L4").Value2 = 26139174
$ws.Range("M4").Value2 = -8480319.5
$ws.Range("N4").Value2 = -26139398
$ws.Range("H6").Value2 = 114.77778
$ws.Range("I6").Value2 = 2.3333333
$ws.Range("J6").Value2 = 339.66666
$ws.Range("K6").Value2 = 6.999999900000001
$ws.Range("L6").Value2 = 1018.99998
$ws.Range("M6").Value2 = 106.0000001
$ws.Range("H21").Value2 = 462.5
$ws.Range("I21").Value2 = 462.5
$ws.Range("J21").Value2 = 0
$ws.Range("K21").Value2 = 1387.5
$ws.Range("L21").Value2 = 0
$ws.Range("M21").Value2 = -1214.5
$ws.Range("H22").Value2 = 11905.667
$ws.Range("I22").Value2 = 646.6667
$ws.Range("J22").Value2 = 17535.166
$ws.Range("K22").Value2 = 1940.0001
$ws.Range("L22").Value2 = 52605.49800000001
$ws.Range("M22").Value2 = -1771.0001
$ws.Range("N22").Value2 = -52943.49800000001
$ws.Range("H27").Value2 = 11905.667
$ws.Range("I27").Value2 = 646.6667
$ws.Range("J27").Value2 = 17535.166
$ws.Range("K27").Value2 = 1940.0001
$ws.Range("L27").Value2 = 52605.49800000001
$ws.Range("M27").Value2 = -1838.0001
$ws.Range("N27").Value2 = -52809.49800000001
$ws.Range("H38").Value2 = 76923110
$ws.Range("I38").Value2 = 250000050
$ws.Range("J38").Value2 = 22.444445
$ws.Range("K38").Value2 = 750000150
$ws.Range("L38").Value2 = 67.33333500000001
$ws.Range("M38").Value2 = -749999803
$ws.Range("N38").Value2 = -761.333335
$ws.Range("H70").Value2 = 499.5
$ws.Range("I70").Value2 = 499.5
$ws.Range("J70").Value2 = 0
$ws.Range("K70").Value2 = 1498.5
$ws.Range("L70").Value2 = 0
$ws.Range("M70").Value2 = -1183.5
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value2 = 499.5
$ws.Range("I73").Value2 = 499.5
$ws.Range("J73").Value2 = 0
$ws.Range("K73").Value2 = 1498.5
$ws.Range("L73").Value2 = 0
$ws.Range("M73").Value2 = -406.5
$ws.Range("N73").ClearContents()
$ws.Range("H74").Value2 = 10000
$ws.Range("I74").Value2 = 10000
$ws.Range("J74").Value2 = 0
$ws.Range("K74").Value2 = 30000
$ws.Range("L74").Value2 = 0
$ws.Range("M74").Value2 = -28939
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value2 = 10000
$ws.Range("I77").Value2 = 10000
$ws.Range("J77").Value2 = 0
$ws.Range("K77").Value2 = 90000
$ws.Range("L77").Value2 = 0
$ws.Range("M77").Value2 = -84696
$ws.Range("N77").ClearContents()
$ws.Range("H98").Value2 = 529.5
$ws.Range("I98").Value2 = 188.33333
$ws.Range("J98").Value2 = 622.5454999999999
$ws.Range("K98").Value2 = 564.99999
$ws.Range("L98").Value2 = 1867.6365
$ws.Range("M98").Value2 = 933.00001
$ws.Range("N98").Value2 = -4863.6365
$ws.Range("H100").Value2 = 0
$ws.Range("I100").Value2 = 0
$ws.Range("J100").Value2 = 0
$ws.Range("K100").Value2 = 0
$ws.Range("L100").Value2 = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H107").Value2 = 1677
$ws.Range("I107").Value2 = 2250
$ws.Range("J107").Value2 = 1562.4
$ws.Range("K107").Value2 = 6750
$ws.Range("L107").Value2 = 4687.200000000001
$ws.Range("M107").Value2 = -4830
$ws.Range("N107").Value2 = -8527.200000000001
$ws.Range("H109").Value2 = 4936.8335
$ws.Range("I109").Value2 = 4936.8335
$ws.Range("J109").Value2 = 0
$ws.Range("K109").Value2 = 14810.5005
$ws.Range("L109").Value2 = 0
$ws.Range("M109").Value2 = -13770.5005
$ws.Range("H117").Value2 = 3221.2
$ws.Range("I117").Value2 = 864.5
$ws.Range("J117").Value2 = 4792.3335
$ws.Range("K117").Value2 = 2593.5
$ws.Range("L117").Value2 = 14377.0005
$ws.Range("M117").Value2 = 848.5
$ws.Range("N117").Value2 = -21261.0005
$ws.Range("H121").Value2 = 14494547
$ws.Range("I121").Value2 = 2450.5
$ws.Range("J121").Value2 = 19609404
$ws.Range("K121").Value2 = 7351.5
$ws.Range("L121").Value2 = 58828212
$ws.Range("M121").Value2 = -6041.5
$ws.Range("N121").Value2 = -58830832
$ws.Range("H131").Value2 = 12827539
$ws.Range("I131").Value2 = 33334386
$ws.Range("J131").Value2 = 10759.5625
$ws.Range("K131").Value2 = 100003158
$ws.Range("L131").Value2 = 32278.6875
$ws.Range("M131").Value2 = -99998118
$ws.Range("N131").Value2 = -42358.6875
$ws.Range("H139").Value2 = 29413980
$ws.Range("I139").Value2 = 31252134
$ws.Range("J139").Value2 = 3500
$ws.Range("K139").Value2 = 93756402
$ws.Range("L139").Value2 = 10500
$ws.Range("M139").Value2 = -93751262
$ws.Range("H140").Value2 = 19232596
$ws.Range("I140").Value2 = 35715476
$ws.Range("J140").Value2 = 2568
$ws.Range("K140").Value2 = 107146428
$ws.Range("L140").Value2 = 7704
$ws.Range("M140").Value2 = -107141248
$ws.Range("N140").Value2 = -18064

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 84312.086
$ws.Range("I2").Value2 = 1171.5555
$ws.Range("J2").Value2 = 333733.66
$ws.Range("K2").Value2 = 1171.5555
$ws.Range("L2").Value2 = 333733.66
$ws.Range("M2").Value2 = -1058.5555
$ws.Range("H70").Value2 = 8353.462
$ws.Range("I70").Value2 = 7037.375
$ws.Range("J70").Value2 = 10459.2
$ws.Range("K70").Value2 = 7037.375
$ws.Range("L70").Value2 = 10459.2
$ws.Range("M70").Value2 = -6767.375
$ws.Range("N70").Value2 = -10999.2
$ws.Range("H73").Value2 = 8353.462
$ws.Range("I73").Value2 = 7037.375
$ws.Range("J73").Value2 = 10459.2
$ws.Range("K73").Value2 = 7037.375
$ws.Range("L73").Value2 = 10459.2
$ws.Range("M73").Value2 = -6101.375
$ws.Range("N73").Value2 = -12331.2
$ws.Range("H102").Value2 = 1723.5581
$ws.Range("I102").Value2 = 1532.1714
$ws.Range("J102").Value2 = 2560.875
$ws.Range("K102").Value2 = 1532.1714
$ws.Range("L102").Value2 = 2560.875
$ws.Range("M102").Value2 = 89.82860000000005
$ws.Range("N102").Value2 = -5804.875
$ws.Range("H118").Value2 = 50000
$ws.Range("I118").Value2 = 0
$ws.Range("J118").Value2 = 50000
$ws.Range("K118").Value2 = 0
$ws.Range("L118").Value2 = 50000
$ws.Range("N118").Value2 = -53314
$ws.Range("H126").Value2 = 4956
$ws.Range("I126").Value2 = 3941.3333
$ws.Range("J126").Value2 = 8000
$ws.Range("K126").Value2 = 11823.9999
$ws.Range("L126").Value2 = 24000
$ws.Range("M126").Value2 = -9353.999899999999
$ws.Range("H132").Value2 = 6755.25
$ws.Range("I132").Value2 = 4853.95
$ws.Range("J132").Value2 = 9924.083000000001
$ws.Range("K132").Value2 = 14561.85
$ws.Range("L132").Value2 = 29772.249
$ws.Range("M132").Value2 = -12031.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 4317.5557
$ws.Range("I7").Value2 = 4251
$ws.Range("J7").Value2 = 4850
$ws.Range("K7").Value2 = 4251
$ws.Range("L7").Value2 = 4850
$ws.Range("M7").Value2 = -4139
$ws.Range("H9").Value2 = 115.666664
$ws.Range("I9").Value2 = 115.666664
$ws.Range("J9").Value2 = 0
$ws.Range("K9").Value2 = 115.666664
$ws.Range("L9").Value2 = 0
$ws.Range("M9").Value2 = 108.333336
$ws.Range("H16").Value2 = 930.425
$ws.Range("I16").Value2 = 881.56757
$ws.Range("J16").Value2 = 1533
$ws.Range("K16").Value2 = 881.56757
$ws.Range("L16").Value2 = 1533
$ws.Range("M16").Value2 = -711.56757
$ws.Range("N16").Value2 = -1873
$ws.Range("H22").Value2 = 4013.24
$ws.Range("I22").Value2 = 2178.6667
$ws.Range("J22").Value2 = 5045.1875
$ws.Range("K22").Value2 = 2178.6667
$ws.Range("L22").Value2 = 5045.1875
$ws.Range("M22").Value2 = -1883.6667
$ws.Range("N22").Value2 = -5635.1875
$ws.Range("H27").Value2 = 4013.24
$ws.Range("I27").Value2 = 2178.6667
$ws.Range("J27").Value2 = 5045.1875
$ws.Range("K27").Value2 = 2178.6667
$ws.Range("L27").Value2 = 5045.1875
$ws.Range("M27").Value2 = -2071.6667
$ws.Range("N27").Value2 = -5259.1875
$ws.Range("H39").Value2 = 36263.5
$ws.Range("I39").Value2 = 37529.5
$ws.Range("J39").Value2 = 34997.5
$ws.Range("K39").Value2 = 37529.5
$ws.Range("L39").Value2 = 34997.5
$ws.Range("M39").Value2 = -37069.5
$ws.Range("N39").Value2 = -35917.5
$ws.Range("H46").Value2 = 4508.091
$ws.Range("I46").Value2 = 1311.2667
$ws.Range("J46").Value2 = 7172.1113
$ws.Range("K46").Value2 = 1311.2667
$ws.Range("L46").Value2 = 7172.1113
$ws.Range("M46").Value2 = -1123.2667
$ws.Range("N46").Value2 = -7548.1113
$ws.Range("H55").Value2 = 594.7222
$ws.Range("I55").Value2 = 635.3333
$ws.Range("J55").Value2 = 513.5
$ws.Range("K55").Value2 = 635.3333
$ws.Range("L55").Value2 = 513.5
$ws.Range("M55").Value2 = -462.3333
$ws.Range("N55").Value2 = -859.5
$ws.Range("H126").Value2 = 4317.5557
$ws.Range("I126").Value2 = 4251
$ws.Range("J126").Value2 = 4850
$ws.Range("K126").Value2 = 12753
$ws.Range("L126").Value2 = 14550
$ws.Range("M126").Value2 = -10283
$ws.Range("H128").Value2 = 171110
$ws.Range("I128").Value2 = 0
$ws.Range("J128").Value2 = 171110
$ws.Range("K128").Value2 = 0
$ws.Range("L128").Value2 = 171110
$ws.Range("N128").Value2 = -181070
$ws.Range("H132").Value2 = 6767.159
$ws.Range("I132").Value2 = 6165.778
$ws.Range("J132").Value2 = 9473.375
$ws.Range("K132").Value2 = 18497.334
$ws.Range("L132").Value2 = 28420.125
$ws.Range("M132").Value2 = -15967.334
$ws.Range("N132").Value2 = -33480.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 19429.889
$ws.Range("I62").Value2 = 23328.166
$ws.Range("J62").Value2 = 11633.333
$ws.Range("K62").Value2 = 23328.166
$ws.Range("L62").Value2 = 11633.333
$ws.Range("M62").Value2 = -22704.166
$ws.Range("N62").Value2 = -12881.333
$ws.Range("H65").Value2 = 19429.889
$ws.Range("I65").Value2 = 23328.166
$ws.Range("J65").Value2 = 11633.333
$ws.Range("K65").Value2 = 116640.83
$ws.Range("L65").Value2 = 58166.665
$ws.Range("M65").Value2 = -113520.83
$ws.Range("N65").Value2 = -64406.665
$ws.Range("H81").Value2 = 10238.28
$ws.Range("I81").Value2 = 6795.8
$ws.Range("J81").Value2 = 12533.267
$ws.Range("K81").Value2 = 13591.6
$ws.Range("L81").Value2 = 25066.534
$ws.Range("M81").Value2 = -12530.6
$ws.Range("H84").Value2 = 10238.28
$ws.Range("I84").Value2 = 6795.8
$ws.Range("J84").Value2 = 12533.267
$ws.Range("K84").Value2 = 67958
$ws.Range("L84").Value2 = 125332.67
$ws.Range("M84").Value2 = -62654
$ws.Range("H100").Value2 = 3406.2
$ws.Range("I100").Value2 = 1131.3334
$ws.Range("J100").Value2 = 4381.143
$ws.Range("K100").Value2 = 2262.6668
$ws.Range("L100").Value2 = 8762.286
$ws.Range("M100").Value2 = -1721.6668
$ws.Range("N100").Value2 = -9844.286
$ws.Range("H127").Value2 = 36796.668
$ws.Range("I127").Value2 = 60390
$ws.Range("J127").Value2 = 25000
$ws.Range("K127").Value2 = 60390
$ws.Range("L127").Value2 = 25000
$ws.Range("M127").Value2 = -55430
$ws.Range("N127").Value2 = -34920
$ws.Range("H128").Value2 = 0
$ws.Range("I128").Value2 = 0
$ws.Range("J128").Value2 = 0
$ws.Range("K128").Value2 = 0
$ws.Range("L128").Value2 = 0
$ws.Range("M128").ClearContents()
$ws.Range("H132").Value2 = 4727.523
$ws.Range("I132").Value2 = 3611.3794
$ws.Range("J132").Value2 = 6885.4
$ws.Range("K132").Value2 = 10834.1382
$ws.Range("L132").Value2 = 20656.2
$ws.Range("M132").Value2 = -8304.138199999999
$ws.Range("N132").Value2 = -25716.2
$ws.Range("H135").Value2 = 89999.664
$ws.Range("I135").Value2 = 0
$ws.Range("J135").Value2 = 89999.664
$ws.Range("K135").Value2 = 0
$ws.Range("L135").Value2 = 89999.664
$ws.Range("N135").Value2 = -100139.664
$ws.Range("H136").Value2 = 3134.9023
$ws.Range("I136").Value2 = 2260.3704
$ws.Range("J136").Value2 = 4821.5
$ws.Range("K136").Value2 = 6781.111199999999
$ws.Range("L136").Value2 = 14464.5
$ws.Range("M136").Value2 = -4231.111199999999
